$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@(44344, 0, 8, 86.16975441619991)
    ,@(44345, 0, 6, 64.62731581214993)
    ,@(44346, 0, 2, 21.54243860404998)
    ,@(44347, 0, 1, 10.77121930202499)
    ,@(44348, 0, 1, 10.77121930202499)
    ,@(44349, 0, 1, 10.77121930202499)
    ,@(44350, 0, 0, 0)
    ,@(44351, 4, 4, 43.08487720809995)
    ,@(44352, 0, 4, 43.08487720809995)
    ,@(44353, 0, 4, 43.08487720809995)
    ,@(44354, 0, 4, 43.08487720809995)
    ,@(44355, 0, 4, 43.08487720809995)
    ,@(44356, 0, 4, 43.08487720809995)
    ,@(44357, 0, 4, 43.08487720809995)
    ,@(44358, 0, 0, 0)
    ,@(44359, 1, 1, 10.77121930202499)
    ,@(44360, 0, 1, 10.77121930202499)
    ,@(44361, 0, 1, 10.77121930202499)
    ,@(44362, 4, 5, 53.85609651012494)
    ,@(44363, 0, 5, 53.85609651012494)
    ,@(44364, 0, 5, 53.85609651012494)
    ,@(44365, 0, 5, 53.85609651012494)
    ,@(44366, 0, 4, 43.08487720809995)
    ,@(44367, 0, 4, 43.08487720809995)
    ,@(44368, 0, 4, 43.08487720809995)
    ,@(44369, 1, 1, 10.77121930202499)
    ,@(44370, 1, 2, 21.54243860404998)
    ,@(44371, 0, 2, 21.54243860404998)
    ,@(44372, 0, 2, 21.54243860404998)
    ,@(44373, 0, 2, 21.54243860404998)
    ,@(44374, 1, 3, 32.31365790607497)
    ,@(44375, 0, 3, 32.31365790607497)
)

$startRow = 270
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}

# Match the date-column style used by the rest of column A (border, bold,
# centered, custom date/time number format) by copying the format of the
# last existing date cell down into the newly added A270:A301 range.
$ws.Range("A269").Copy()
$ws.Range("A270:A301").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A1").Select()
